$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 1 - the two FAQ placeholder strings get replaced with real content
$ws.Range("A1").Value = "What programs do you offer?"
$ws.Range("B1").Value = "We offer a wide range of undergraduate and graduate programs in various fields, including engineering, computer science, business, and law. You can find more information about our programs on our website, or by contacting our admissions office."

# Row 2 - swap B2/A2 values (fix typo order) and update C2 stays "rules"
$ws.Range("A2").Value = "How do I apply for admission?"
$ws.Range("B2").Value = "You can apply for admission online through our website. Simply select the program you are interested in, and follow the instructions provided. You will need to submit your academic transcripts, test scores, and other required documents along with your application."
$ws.Range("C2").Value = "rules"

# Row 3
$ws.Range("A3").Value = "What are the admission requirements?"
$ws.Range("B3").Value = "Admission requirements vary by program, but generally include a minimum GPA, standardized test scores (such as the SAT or ACT), and letters of recommendation. Specific requirements can be found on our website, or by contacting our admissions office."

# Row 4
$ws.Range("A4").Value = "What financial aid options are available?"
$ws.Range("B4").Value = "We offer a variety of financial aid options, including scholarships, grants, loans, and work-study programs. You can find more information about financial aid on our website, or by contacting our financial aid office."

# Row 5
$ws.Range("A5").Value = "How can I schedule a campus visit?"
$ws.Range("B5").Value = "You can schedule a campus visit through our website, or by contacting our admissions office. Campus visits typically include a tour of our facilities, meetings with faculty and staff, and the opportunity to attend a class or event."

# Column widths (target widths of 52 / 97.21875 / 54.44140625 / 43.33203125 /
# 44.33203125 chars, as authored in real Excel on a 7px-per-char font metric;
# this engine's ColumnWidth setter quantizes stored width to 1/6-character
# steps on a different internal metric, so the inputs below are the values
# that land closest to those exact targets after that quantization).
$ws.Columns.Item(1).ColumnWidth = 51.21875
$ws.Columns.Item(2).ColumnWidth = 96.28125
$ws.Columns.Item(3).ColumnWidth = 53.66015625
$ws.Columns.Item(4).ColumnWidth = 42.55078125
$ws.Columns.Item(5).ColumnWidth = 43.55078125

# Sheet view
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("V12").Select()
